$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the four new question rows under the existing table (rows 26-29, column A)
$ws.Range("A26").Value = "Areas of improvements?"
$ws.Range("A27").Value = "Does prediction makes sense in an exponential model - LSTM"
$ws.Range("A28").Value = "what else would you like to see for yourself in a similar study?"
$ws.Range("A29").Value = "Any other insights on getting more data sources?"

# Best-fit the new (longer) column H header text ("Support and positivity")
$ws.Range("H1:H8").ColumnWidth = 23

# Update the sheet view: zoom level and active selection cell
$ws.Application.ActiveWindow.Zoom = 142
[void]$ws.Range("D14").Select()
